$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column D (current "Terms Typically Offered")
# so that D, E, F become new blank columns and the old D shifts to G.
$ws.Range("D:F").Insert()

# New header row values
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# New column values for data rows 2-9
$ws.Range("D2:F9").Value = "NA"
